$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) from 45170 to 45174
# for data rows 2 through 13.
for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 3).Value = 45174
}
